# Finnhub RTD demo workbook update
# - Reword the "paste your token" helper cell (B1) to the new API-key message
# - Move the active selection from F16 to B2
# - Re-apply the RTD() formulas for the quote grid (B4:D7) so their values
#   get refreshed (the RTD provider's "getting data" placeholder now
#   resolves/ defaults to 0 instead of the stale cached numbers/timestamps
#   that were baked into the workbook the last time a live RTD server was
#   attached).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1: update the instructional formula/text -----------------------------
$ws.Range("B1").Formula = '=IF(A1="","<--Copy paste your Finnhub API key to A1 cell. Go to https://finnhub.io/dashboard to get your API key.",RTD("finnhub","","set_token",A1))'

# --- Quote grid: refresh the RTD() formulas for rows 4-7 --------------------
$rows = 4,5,6,7
foreach ($r in $rows) {
    $ws.Range("B$r").Formula = '=RTD("finnhub","",A' + $r + ',"last_price")'
    $ws.Range("C$r").Formula = '=RTD("finnhub","",A' + $r + ',"volume")'
    $ws.Range("D$r").Formula = '=RTD("finnhub","",A' + $r + ',"last_update_time")'
}

# --- Selection: move from F16 to B2 -----------------------------------------
[void]$ws.Range("B2").Select()
